$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

# Regular value updates (price and volume columns)
Set-TextValue "D2" "28.818.69"
Set-TextValue "E2" "  +7.72%  "
Set-TextValue "D3" "1.812.75"
Set-TextValue "E3" "  +5.07%  "
Set-TextValue "D4" "0.9996"
Set-TextValue "E4" "  +0.22%  "
Set-TextValue "D5" "249.87"
Set-TextValue "E5" "  +3.70%  "
Set-TextValue "D6" "0.9997"
Set-TextValue "E6" "  +0.18%  "
Set-TextValue "D7" "0.4956"
Set-TextValue "E7" "  +2.41%  "
Set-TextValue "D8" "0.2789"
Set-TextValue "E8" "  +8.05%  "
Set-TextValue "D9" "0.06406"
Set-TextValue "E9" "  +3.39%  "
Set-TextValue "D10" "1.818.33"
Set-TextValue "E10" "  +5.30%  "
Set-TextValue "E11" "  +5.06%  "
Set-TextValue "D12" "0.07139"
Set-TextValue "E12" "  +3.40%  "
Set-TextValue "D13" "0.6507"
Set-TextValue "E13" "  +7.21%  "
Set-TextValue "D14" "83.97"
Set-TextValue "E14" "  +9.20%  "
Set-TextValue "E15" "  +5.17%  "
Set-TextValue "D16" "28.801.89"
Set-TextValue "E16" "  +8.45%  "
Set-TextValue "D17" "0.9999"
Set-TextValue "E17" "  +0.16%  "
Set-TextValue "D18" "0.000007408"
Set-TextValue "E18" "  +3.58%  "
Set-TextValue "D19" "0.9992"
Set-TextValue "E19" "  +0.19%  "
Set-TextValue "E20" "  +7.06%  "
Set-TextValue "D21" "2.044.98"
Set-TextValue "E21" "  +4.89%  "
Set-TextValue "D22" "4.610"
Set-TextValue "E22" "  +4.18%  "
Set-TextValue "D23" "8.902"
Set-TextValue "E23" "  +3.88%  "
Set-TextValue "E24" "  +5.74%  "
Set-TextValue "D25" "143.19"
Set-TextValue "E25" "  +4.50%  "
Set-TextValue "D26" "133.00"
Set-TextValue "E26" "  +25.58%  "
Set-TextValue "D27" "16.58"
Set-TextValue "E27" "  +8.89%  "
Set-TextValue "E28" "  +6.93%  "
Set-TextValue "D29" "1.396"
Set-TextValue "E29" "  +1.27%  "
Set-TextValue "D30" "4.172"
Set-TextValue "E30" "  +5.96%  "
Set-TextValue "D31" "0.08367"
Set-TextValue "E31" "  +5.27%  "
Set-TextValue "D32" "3.857"
Set-TextValue "E32" "  +4.41%  "
Set-TextValue "D33" "0.04951"
Set-TextValue "E33" "  +10.49%  "
Set-TextValue "D34" "1.092"
Set-TextValue "E34" "  +8.41%  "
Set-TextValue "D39" "0.9578"
Set-TextValue "E39" "  +3.48%  "
Set-TextValue "D40" "6.083"
Set-TextValue "E40" "  +7.60%  "
Set-TextValue "D41" "0.01594"
Set-TextValue "E41" "  +6.78%  "
Set-TextValue "D42" "0.9995"
Set-TextValue "E42" "  +0.22%  "
Set-TextValue "D43" "100.74"
Set-TextValue "E43" "  +1.14%  "
Set-TextValue "D44" "0.4105"
Set-TextValue "E44" "  +7.08%  "
Set-TextValue "D45" "7.224"
Set-TextValue "E45" "  +5.46%  "
Set-TextValue "D46" "0.1226"
Set-TextValue "E46" "  +6.07%  "
Set-TextValue "D47" "0.05513"
Set-TextValue "D48" "8.166"
Set-TextValue "E48" "  +3.44%  "
Set-TextValue "E49" "  +5.40%  "
Set-TextValue "D50" "1.313"
Set-TextValue "E50" "  +6.93%  "
Set-TextValue "D51" "0.3635"
Set-TextValue "E51" "  +8.18%  "
Set-TextValue "B35" "ImmutableX"
Set-TextValue "C35" "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextValue "D35" "0.6802"
Set-TextValue "E35" "  +9.62%  "
Set-TextValue "B36" "HuobiToken"
Set-TextValue "C36" "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
Set-TextValue "D36" "2.710"
Set-TextValue "E36" "  +4.55%  "
Set-TextValue "B37" "RenderToken"
Set-TextValue "C37" "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue "D37" "2.271"
Set-TextValue "E37" "  +11.88%  "
Set-TextValue "B38" "MXToken"
Set-TextValue "C38" "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextValue "D38" "2.779"
Set-TextValue "E38" "  +13.72%  "
